# Review_402: update title/date, rewrite review body paragraphs for the new
# paper ("Representation Alignment for Generation"), and append the new arxiv link.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $old"
    }
}

# title line 1 (date)
Replace-Text "המאמר היומי של מייק - 18.02.25" "המאמר היומי של מייק - 16.02.25"

# title line 2 (paper name)
Replace-Text "THINKING LLMS: GENERAL INSTRUCTION FOLLOWING WITH THOUGHT GENERATION" "Representation Alignment for Generation: Training Diffusion Transformers is Easier than you Think"

# paragraph 2 (intro)
Replace-Text "סקירה מספר 400 - כדי לא להכביד עליכם יותר מדי בחרתי מאמר קליל יחסית והסקירה הולכת להיות בלי נוסחאות ודי קצרה. המאמר מציע שיטה קצת במהות דומה Group Relative Preference Optimization או GRPO בקצרה שעשתה הרבה כותרות לאחרונה. ותיכף אני הולך להסביר למה אני מתכוון כאן. רק אציין שהמאמר מציע שיטה להגברת יכולת הנמקה כללית של מודל ולא מתמקד רק בשאלות תכנות ובעיות מתמטיות." "לוקחים פסק זמן קטן מ-LLMs וסוקרים מאמר על מודלי דיפוזיה גנרטיביים. המאמר מציע שיטה די אינטואיטיבית לשיפור  ביצועים של מודלים אלו על ידי הוספת איבר רגולריזציה ה״מיישר״ את הייצוגים הפנימיים של המודל עם אלו של אנקודרים חזקים כמו DiNOV2. יישור זה משפר את איכות התמונות שהמודל מגנרט."

# paragraph 3 (diffusion background)
Replace-Text "המאמר מציע שיטת טיוב (fine-tune) למודלי שפה המתקדמת בהקנייתם יכולת הנמקה (reasoning) למודלי שפה ללא צורך בדאטה מתויג. המאמר מציע לבצע אימון בסגנון RLHF אבל להבדיל מהדיפסיק (הממציאה של GRPO), המחברים הציעו להשתמש בשיטת DPO שלא משתמשת בפונקציית התגמול כלל. אציין ש-GRPO לא מאמנת מודל תגמול (reward) כמו ש-PPO עושה אלא משתמשת בנכונות התשובה והפורמט שלה כפונקצית תגמול." "נתחיל מרקע קצרצר על מודלי דיפוזיה גנרטיביים. מודלים אלו מאומנים לגנרט תמונות (למשל בהינתן תיאור טקסטואלי) על ידי הסרה הדרגתית של הרעש. המודל מתחיל מרעש טהור (בד״כ גאוסי) ולאט לאט הופכים אותו לתמונה (או פיסת דאטה מדומיין אחר). המודל מאומן על תמונות מורעשות עם רמות שונות של רעש(=איטרציות) כאשר באימון המודל לומד להסיר כמות קטנה של רעש (מאיטרציה t לאיטרציה t -1). בחירה של הייפר-הפרמטרים  של תהליך ההרעשה היא מרכיב קריטי לאיכות גנרוט של המודל המאומן."

# paragraph 4 (probability flow)
Replace-Text "אז מה משותף בין GRPO לבין השיטה המוצעת במאמר? שניהם למעשה מציעים לא לקנוס את המודל על תהליך החשיבה (שעלול להיות לא נכון אך להוביל לתשובה הנכונה) אלא לשפוט אותו רק על בסיס נכונות התשובה של המודל (כאמור GRPO גם קונס על אי עמידה בפורמט של התשובה). אחרי שהבנו את הקשרים המהותיים של השיטה המוצעות עם השיטות המפורסמות בואו נצלול למה שהמאמר מציע." "תהליך זה(הרעשה) ניתן לתאר באמצעות משוואות דיפרנציאלית של זרימה הסתברותית (probability flow) המתאר השתנות (גרדיאנט) הדאטה המורעש עם קצב/מהירות הרעשה (velocity) שנסמן אותו  (הפתרון של משוואה זו מתפלג לפי ההתפלגות של הדאטה המורעש). קצב הרעשה ניתן לשערך עם המודל (=רשת) בהתבסס על דגימות הדאטה המורעש ו-. לאחר מכן ניתן לפתור את משוואות הזרימה ההסתברותית עם השערוך של  (בכיוון ההפוך - כלומר החל מרעש טהור) עם שיטת איולר למשל. שיטות אלו נקראות stochastic interpoland. נציין שיש שיטות המבוססות על פתרון נומרי של משוואה דיפרנצאלית סטוכסטית שמתארת את השתנות הדאטה כפונקציה של פונקציית score שהיא לוגריתם של פונקציית התפלגות של דאטה מורעש."

# paragraph 5 (latent diffusion)
Replace-Text "כאמור המאמר מציע לטייב יכולת הנמקה של מודל שפה ללא שימוש בדאטה מתיוג עם RLHF. כמו שאתם זוכרים RLHF עם DPO דורש זוגות של תשובות מועדפות ופחות מועדפות. מכיוון שאמרנו שהשיטה לא דורשת דאטה מתויג אז אתם יכולים לנחש שבניית הזוגות נעשית על ידי מודל שפה שופט שבוחר תשובות טובות ורעות בדומה לשיטת RLAIF שזה קיצור של Reinforcement Learning from AI Feedback. מודל השופט מופעל על תשובות (ולא שרשרת הנמקה!) של המודל המאומן ומחליט מה בין תשובות היא הטובה והגרועה ביותר. זוגות אלו משמשים לאימון המודל בצורת DPO. כמובן שיש פה גם הנדסת של מטה-פרומפט הגורם למודל ״לחשוב״ אבל שרשרת חשיבה זו לא משתתפת באימון המודל." "אוקיי, אחרי הסיבוך הזה החיים נהיים קצת יותר קלים. מודלי דיפוזיה היום הם לרוב מודלים לטנטים כאשר הגנרוט מתרחש במרחב הייצוג של הדאטה. כלומר המודל מאומן לשחזר ייצוג לטנטי מרעש ואז מפעילים את הדקודר כדי לבנות תמונה מהייצוג המשוחזר. הייצוג של התמונה ההתחלתית נוצר על ידי האנקודר. המחברים טוענים שהייצוגים הלטנטיים המורעשים אינם ״חזקים מספיק״ כלומר פחות משקפים את האספקטים הסמנטיים של התמונה. "

# paragraph 6 -> becomes proposed method text (was the old arxiv link paragraph)
Replace-Text "https://arxiv.org/abs/2410.10630" "המחברים מציעים להעשיר את הייצוגים האלו על ידי הוספה של איבר רגולריזציה שמטרתו לקרב ייצוגים אלה (של התמונת המרועשות) לייצוג המופק על ידי אנקודר חזק (כמו DINOV2). לוס זה מתווסף ללוס הרגיל של מודל דיפוזיה ונטען במאמר שזה משפר את איכות התמונות המגונרטות וגם תורם ליציבות האימון."

# Append a new final paragraph holding the new arxiv link
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "https://arxiv.org/abs/2410.06940"

Write-Output ("Done. Final paragraph count: " + $d.Paragraphs.Count)
